$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B2 value from 2 to 1
$ws.Range("B2").Value = 1

# Add new row 3 data
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1

# Copy style from A2 to A3 so it matches (bold font, border, centered alignment)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Ensure value stays correct after paste special (paste formats only, but just in case set again)
$ws.Range("A3").Value = 1
